# Updated symbol list on Sun Jan  1 22:33:02 UTC 2023 with GitHub Actions
#
# Refresh the cryptocurrency "Price" (column D) and "Volume(1h)" (column E)
# figures on the active sheet to the latest scraped values.
#
# The sheet stores these figures as literal text (e.g. "244.82", "-0.87%"),
# not as numbers, so each value is written with a leading apostrophe to stop
# Excel from auto-converting number-/percent-looking text into a numeric
# value. The apostrophe prefix marks the cell with a "quote prefix" style;
# resetting the cell style back to "Normal" afterwards clears that marker
# again while leaving the text value (and its Text cell type) intact, so
# the cell ends up formatted exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "244.78" },
    @{ Cell = "E2"; Value = "-0.86%" },
    @{ Cell = "D3"; Value = "27.43" },
    @{ Cell = "E3"; Value = "4.97%" },
    @{ Cell = "D4"; Value = "5.116" },
    @{ Cell = "E4"; Value = "0.58%" },
    @{ Cell = "D5"; Value = "0.05682" },
    @{ Cell = "E5"; Value = "1.51%" },
    @{ Cell = "D6"; Value = "6.521" },
    @{ Cell = "E6"; Value = "0.66%" },
    @{ Cell = "D7"; Value = "0.8199" },
    @{ Cell = "E7"; Value = "0.67%" },
    @{ Cell = "D8"; Value = "0.8522" },
    @{ Cell = "E8"; Value = "0.99%" },
    @{ Cell = "D9"; Value = "0.06951" },
    @{ Cell = "E9"; Value = "-0.48%" },
    @{ Cell = "D10"; Value = "0.02876" },
    @{ Cell = "E10"; Value = "2.07%" },
    @{ Cell = "E11"; Value = "0.13%" },
    @{ Cell = "D12"; Value = "0.001516" },
    @{ Cell = "E12"; Value = "-0.10%" },
    @{ Cell = "D13"; Value = "0.04028" },
    @{ Cell = "E13"; Value = "-13.47%" },
    @{ Cell = "D14"; Value = "0.0005979" },
    @{ Cell = "E14"; Value = "0.31%" },
    @{ Cell = "D15"; Value = "0.006215" },
    @{ Cell = "E15"; Value = "0.57%" },
    @{ Cell = "D16"; Value = "3.513" },
    @{ Cell = "E16"; Value = "-2.65%" },
    @{ Cell = "D17"; Value = "3.009" },
    @{ Cell = "E17"; Value = "-0.37%" },
    @{ Cell = "D18"; Value = "2.231" },
    @{ Cell = "E18"; Value = "8.55%" },
    @{ Cell = "D19"; Value = "0.3189" },
    @{ Cell = "E19"; Value = "2.46%" },
    @{ Cell = "D20"; Value = "0.1333" },
    @{ Cell = "E20"; Value = "-0.01%" },
    @{ Cell = "D21"; Value = "0.03208" },
    @{ Cell = "E21"; Value = "-0.04%" },
    @{ Cell = "E22"; Value = "-0.03%" },
    @{ Cell = "D23"; Value = "3.558" },
    @{ Cell = "E23"; Value = "-4.92%" },
    @{ Cell = "E25"; Value = "-2.70%" },
    @{ Cell = "E26"; Value = "-2.03%" },
    @{ Cell = "E27"; Value = "22.86%" },
    @{ Cell = "D28"; Value = "0.0001406" },
    @{ Cell = "E28"; Value = "-27.43%" },
    @{ Cell = "E40"; Value = "1.66%" },
    @{ Cell = "D41"; Value = "0.005990" },
    @{ Cell = "E41"; Value = "-2.64%" },
    @{ Cell = "D42"; Value = "0.1060" },
    @{ Cell = "E42"; Value = "0.55%" },
    @{ Cell = "D43"; Value = "0.002350" },
    @{ Cell = "E43"; Value = "-10.21%" },
    @{ Cell = "D44"; Value = "0.009715" },
    @{ Cell = "E44"; Value = "20.50%" },
    @{ Cell = "D45"; Value = "0.00005122" },
    @{ Cell = "E45"; Value = "-5.07%" },
    @{ Cell = "E46"; Value = "-0.04%" },
    @{ Cell = "E47"; Value = "-30.35%" },
    @{ Cell = "D48"; Value = "0.002516" },
    @{ Cell = "E48"; Value = "3.69%" },
    @{ Cell = "E49"; Value = "-0.04%" },
    @{ Cell = "E50"; Value = "-0.04%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}
